$wb = $excel.ActiveWorkbook

# --- Update Metadata sheet (Version and Date) ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "0.9.20"
$meta.Range("B8").Value = "2025-11-18T19:57:13-03:00"

# --- Update Elements sheet (Constraint(s) column for Extension.value[x] row) ---
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AJ6").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`npositive-integer:Group size must be a positive integer (greater than 0) {`$this > 0}"
